# Applies the changes described in the commit diff:
#  - rename "StrategySpaceInReference5" -> "StrategySpaceInReference 10"
#  - rename "StrategySpaceInReference7" -> "StrategySpaceInReference 11"
#  - move the active tab / selection from the 1st sheet (AgentSettings) to the
#    3rd sheet (StrategySpaceInReference 11), updating its selected cell to D23
#  - (the x15ac:absPath drive letter and the xr:revisionPtr documentId are
#    machine/session generated values outside the reach of the Excel object
#    model exposed here, so they are left to the runtime)

$wb = $excel.ActiveWorkbook

$wsAgentSettings = $wb.Worksheets.Item(1)
$wsRef10 = $wb.Worksheets.Item(2)
$wsRef11 = $wb.Worksheets.Item(3)

# Rename the strategy-space sheets.
$wsRef10.Name = "StrategySpaceInReference 10"
$wsRef11.Name = "StrategySpaceInReference 11"

# Move selection within AgentSettings doesn't change (stays C25), but it is
# no longer the active/tabSelected sheet.
$wsAgentSettings.Range("C25").Select()

# Make the renamed 3rd sheet the active tab and update its selection.
$wsRef11.Select()
$wsRef11.Range("D23").Select()
